# Fill in the previously-empty "Task 4" confusion-matrix block on the
# "FS-IF-IA" worksheet with the experiment results (t-SNE run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FS-IF-IA")

# Switch to manual calculation so the dependent formulas (F24:F26, M24:M26)
# keep their stale cached "#DIV/0!" results, matching the original author's
# edit (which did not trigger a recalculation).
$excel.Calculation = -4135

# --- Left block (Predicted, "Success rate (%)") ---
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 31.57894736842105

$ws.Range("C25").Value = 0.25
$ws.Range("D25").Value = 0.625
$ws.Range("E25").Value = 0.125

$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.6666666666666666

# --- Right block (Predicted, "Overall cross-validation error (%)") ---
$ws.Range("J24").Value = 0.9777777777777777
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 3.508771929824561

$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = 0.9722222222222222
$ws.Range("L25").Value = 0.027777777777777776

$ws.Range("J26").Value = 0.018518518518518517
$ws.Range("K26").Value = 0.037037037037037035
$ws.Range("L26").Value = 0.9444444444444444
